# The deck originally carries the "Integral" theme (ppt/theme/theme2.xml,
# used by the slide master / all slides) and an orphaned default
# "Office Theme" (ppt/theme/theme1.xml, used only by the notes master).
# The commit swaps those two themes' contents, so the presentation itself
# (slides/slide master) ends up on the plain "Office Theme" colour scheme.
#
# Reproduce that by pushing the stock Office theme colours into the
# presentation's active theme colour scheme (msoThemeColor indices
# 1..12 = dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$scheme.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$scheme.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$scheme.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$scheme.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$scheme.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$scheme.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$scheme.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$scheme.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$scheme.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$scheme.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$scheme.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$scheme.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
